# Correcting market share tab for updated scenario 3s
#
# The "New Product A" market share (row 2) and "Old Product B (SOC)" market
# share (row 3) on the MarketShare sheet had the wrong years' worth of 100%
# market-share flags. New Product A should hold 100% share from 2026 (col L)
# through 2040 (col Z), while Old Product B (SOC) should keep its 100% share
# only from 2018 (col D) through 2025 (col K).

$wb        = $excel.ActiveWorkbook
$wsMarket  = $wb.Worksheets.Item("MarketShare")

# Columns L..Z (2026-2040) move from the "Old Product B (SOC)" row (row 3)
# to the "New Product A" row (row 2); the old row keeps only D..K (2018-2025).
$colsToMove = @("L","M","N","O","P","Q","R","S","T","U","V","W","X","Y","Z")
foreach ($col in $colsToMove) {
    $newProductCell = $col + "2"
    $oldProductCell = $col + "3"
    $wsMarket.Range($newProductCell).Value = 1
    $wsMarket.Range($oldProductCell).Value = $null
}

# MarketShare becomes the active/selected sheet (was "Platform Coverage").
$wsMarket.Activate()

# Update the view's scroll position / selection on the MarketShare sheet.
$excel.ActiveWindow.ScrollColumn = 6
$wsMarket.Range("Y6").Select() | Out-Null
